# Regenerate save_data to use K (strikeouts) instead of Strike# in column G.
# This writes the corrected per-appearance strikeout ("K") totals for
# Jake Diekman's 2023 save_data sheet (column G, rows 2-73).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 2
    3  = 0
    4  = 1
    5  = 2
    6  = 1
    7  = 0
    8  = 3
    9  = 1
    10 = 2
    11 = 2
    12 = 1
    13 = 2
    14 = 2
    15 = 0
    16 = 1
    17 = 0
    18 = 1
    19 = 1
    20 = 0
    21 = 1
    22 = 1
    23 = 1
    24 = 2
    25 = 0
    26 = 0
    27 = 1
    28 = 0
    29 = 1
    30 = 0
    31 = 1
    32 = 1
    33 = 1
    34 = 3
    35 = 3
    36 = 5
    37 = 2
    38 = 2
    39 = 0
    40 = 1
    41 = 1
    42 = 2
    43 = 1
    44 = 0
    45 = 0
    46 = 1
    47 = 1
    48 = 1
    49 = 0
    50 = 1
    51 = 0
    52 = 1
    53 = 0
    54 = 2
    55 = 1
    56 = 2
    57 = 0
    58 = 1
    59 = 0
    60 = 0
    61 = 1
    62 = 0
    63 = 1
    64 = 4
    65 = 0
    66 = 1
    67 = 0
    68 = 0
    69 = 0
    70 = 1
    71 = 1
    72 = 1
    73 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
